# Edit script: update blog visitor-stats worksheet for August 2024 data,
# matching the commit's replacement of the February 2024 table with a new
# August 2024 table (31 days instead of 29, hence 2 extra rows inserted
# before the trailing blank/footer row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header cells (plain text, safe to set directly) ----
$ws.Range("B4").Value = "2024-08-01 ~ 2024-08-31"
$ws.Range("B5").Value = "2024년 09월 04일 19시 05분 23초"

# ---- Make room: August has 31 days vs February's 29, so insert two more
#      data rows above the trailing blank/footer row (old row 37 -> row 39). ----
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()

# ---- Column A: dates (text look like numbers to Excel's type-sniffer, so
#      force the destination range to Text format before the bulk write,
#      then restore the normal format via Paste Special so the XML keeps
#      using the original shared cell style rather than a new "@" style). ----
$dates = @("2024-08-31","2024-08-30","2024-08-29","2024-08-28","2024-08-27","2024-08-26","2024-08-25","2024-08-24","2024-08-23","2024-08-22","2024-08-21","2024-08-20","2024-08-19","2024-08-18","2024-08-17","2024-08-16","2024-08-15","2024-08-14","2024-08-13","2024-08-12","2024-08-11","2024-08-10","2024-08-09","2024-08-08","2024-08-07","2024-08-06","2024-08-05","2024-08-04","2024-08-03","2024-08-02","2024-08-01")
$rngA = $ws.Range("A8:A38")
$rngA.NumberFormat = "@"
$arrA = New-Object 'object[,]' 31,1
for ($i = 0; $i -lt 31; $i++) { $arrA[$i,0] = $dates[$i] }
$rngA.Value = $arrA
$ws.Range("A8").Copy()
$rngA.PasteSpecial(-4122)

# ---- Column B: day-of-week (Korean glyphs, never numeric-looking, so a
#      plain bulk assignment keeps its original text type/style). ----
$weekdays = @("토","금","목","수","화","월","일","토","금","목","수","화","월","일","토","금","목","수","화","월","일","토","금","목","수","화","월","일","토","금","목")
$rngB = $ws.Range("B8:B38")
$arrB = New-Object 'object[,]' 31,1
for ($i = 0; $i -lt 31; $i++) { $arrB[$i,0] = $weekdays[$i] }
$rngB.Value = $arrB

# ---- Columns C:F: numeric-looking counts, need the same Text-forcing
#      treatment as column A. ----
$cdef = @(
    @("25","1","0","24"),
    @("16","0","0","16"),
    @("20","0","0","20"),
    @("21","0","0","21"),
    @("21","0","0","21"),
    @("21","0","0","21"),
    @("12","0","0","12"),
    @("16","0","0","16"),
    @("18","0","0","18"),
    @("12","0","0","12"),
    @("21","1","0","21"),
    @("19","0","0","19"),
    @("17","1","0","16"),
    @("16","0","0","16"),
    @("10","1","0","9"),
    @("21","0","0","21"),
    @("19","1","0","18"),
    @("18","1","0","17"),
    @("19","0","0","19"),
    @("16","0","0","16"),
    @("17","0","0","17"),
    @("14","0","0","14"),
    @("19","0","0","19"),
    @("16","0","0","16"),
    @("6","1","0","5"),
    @("14","1","0","13"),
    @("10","0","0","10"),
    @("17","1","0","16"),
    @("11","0","0","11"),
    @("17","0","0","17"),
    @("14","0","0","14")
)
$rngCF = $ws.Range("C8:F38")
$rngCF.NumberFormat = "@"
$arrCF = New-Object 'object[,]' 31,4
for ($i = 0; $i -lt 31; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arrCF[$i,$j] = $cdef[$i][$j]
    }
}
$rngCF.Value = $arrCF

# ---- Restore the original alternating row-banding style (style index 8 on
#      even data rows, 9 on odd data rows) that the Text-format step above
#      replaced with a freshly minted "@" style. ----
for ($r = 8; $r -le 38; $r++) {
    if ($r % 2 -eq 0) {
        $ws.Range("B9").Copy()
    } else {
        $ws.Range("B8").Copy()
    }
    $ws.Range("C" + $r + ":F" + $r).PasteSpecial(-4122)
}

Write-Output "done"
